$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows just above the old blank separator row (old row 61),
# pushing the separator (-> 63) and summary (-> 64) rows down, and the
# trailing spacer row (old 65 -> 67).
[void]$ws.Rows.Item(61).Insert()
[void]$ws.Rows.Item(61).Insert()

# The insert operation leaves the new rows with formatting copied from the
# row above but no explicit values; wipe them so we can rebuild precisely.
[void]$ws.Range("A61:P62").Clear()

# Re-apply the same per-column formatting used by the data rows above (row 60)
[void]$ws.Range("A60:P60").Copy()
[void]$ws.Range("A61:P61").PasteSpecial(-4122)
[void]$ws.Range("A60:P60").Copy()
[void]$ws.Range("A62:P62").PasteSpecial(-4122)

# --- Row 61: new run at The Spotted Cow, Holbrook on 21/04/2021 ---
$ws.Range("A61").Value = 44307
$ws.Range("B61").Value = "The Spotted Cow"
$ws.Range("C61").Value = "Holbrook"
$ws.Range("D61").Value = "start/end at pub"
$ws.Range("E61").Value = 1.72
$ws.Range("F61").Value = 0.020219907407407409
$ws.Range("G61").Formula = "=F61/E61"
$ws.Range("H61").Value = 1
$ws.Range("I61").Value = 1
$ws.Range("J61").Value = 1
[void]$ws.Range("K61:M61").Clear()
$ws.Range("N61").Value = 1
$ws.Range("O61").Value = "Back after lockdown"
$ws.Range("P61").Formula = "=SUM(H61:N61)*E61"

# --- Row 62: new run at The Spotted Cow, Holbrook on 19/05/2021 ---
$ws.Range("A62").Value = 44335
$ws.Range("B62").Value = "The Spotted Cow"
$ws.Range("C62").Value = "Holbrook"
$ws.Range("D62").Value = "start/end at pub"
$ws.Range("E62").Value = 2.33
$ws.Range("F62").Value = 0.023738425925925923
$ws.Range("G62").Formula = "=F62/E62"
$ws.Range("H62").Value = 1
$ws.Range("I62").Value = 1
[void]$ws.Range("J62:M62").Clear()
$ws.Range("N62").Value = 1
$ws.Range("O62").Value = "A gentle trot to Makeney and back"
$ws.Range("P62").Formula = "=SUM(H62:N62)*E62"

# --- Update the summary row (now row 64) so its ranges include the two new rows ---
$ws.Range("E64").Formula = "=SUM(E5:E62)"
$ws.Range("G64").Formula = "=AVERAGE(G6:G62)"
$ws.Range("H64").Formula = "=SUM(H5:H62)"
$ws.Range("N64").Formula = "=SUM(N5:N62)"
$ws.Range("P64").Formula = "=SUM(P5:P63)"

# Leave the selection on the newly added last row, like the author would have
[void]$ws.Range("A62").Select()
